$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 27, shifting the existing rows 27:40 down to 28:41.
$ws.Rows.Item(27).Insert()

# Populate the new row 27 with the new price-report record.
$ws.Cells.Item(27, 1).Value = 5
$ws.Cells.Item(27, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(27, 3).Value = "Maule"
$ws.Cells.Item(27, 4).Value = 44754
$ws.Cells.Item(27, 5).Value = 7
$ws.Cells.Item(27, 6).Value = 100112043
$ws.Cells.Item(27, 7).Value = "Pepino dulce"
$ws.Cells.Item(27, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 400
$ws.Cells.Item(27, 11).Value = 15000
$ws.Cells.Item(27, 12).Value = 15000
$ws.Cells.Item(27, 13).Value = 15000
$ws.Cells.Item(27, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(27, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(27, 16).Value = 833
$ws.Cells.Item(27, 17).Value = 18
$ws.Cells.Item(27, 18).Value = "Hortaliza"
